# Regenerate save_data to use K (strike count) instead of Strike# (TB),
# recompute std/mean of the dS series, and calculate + write the s_vals
# (simulated strike-touch counts) back into column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values (s_vals), one per data row (rows 2..53), computed from the
# recalculated simulation using dS0/dSF mean & std-dev.
$sVals = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 0
    13 = 2
    14 = 0
    15 = 1
    16 = 0
    17 = 1
    18 = 3
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 2
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 2
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 1
    34 = 1
    35 = 1
    36 = 2
    37 = 0
    38 = 1
    39 = 0
    40 = 1
    41 = 2
    42 = 1
    43 = 1
    44 = 3
    45 = 2
    46 = 2
    47 = 1
    48 = 1
    49 = 1
    50 = 3
    51 = 1
    52 = 2
    53 = 1
}

foreach ($row in $sVals.Keys) {
    $ws.Cells.Item($row, 7).Value = $sVals[$row]
}
